$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.889.29'
$ws.Range("E2").Value = '  +1.54%  '
$ws.Range("D3").Value = '3.413.90'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '578.53'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").Value = '144.22'
$ws.Range("E6").Value = '  +2.47%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("E10").Value = '  +1.00%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '4.001.95'
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("D13").Value = '28.54'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = '3.416.73'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '61.956.13'
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '6.17'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").Value = '14.00'
$ws.Range("E19").Value = '  +3.10%  '
$ws.Range("D20").Value = '9.15'
$ws.Range("E20").Value = '  +3.13%  '
$ws.Range("D21").Value = '391.56'
$ws.Range("E21").Value = '  +2.70%  '
$ws.Range("D22").Value = '74.81'
$ws.Range("E22").Value = '  -1.63%  '
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").Value = '3.558.85'
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = '0.184'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").Value = '7.46'
$ws.Range("E28").Value = '  +3.45%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '8.00'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").Value = '2.14'
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.40'
$ws.Range("E33").Value = '  +3.52%  '
$ws.Range("D34").Value = '23.55'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '5.27'
$ws.Range("E35").Value = '  +6.37%  '
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("D37").Value = '167.66'
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("D38").Value = '1.52'
$ws.Range("E38").Value = '  +4.68%  '
$ws.Range("D39").Value = '3.446.16'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").Value = '28.89'
$ws.Range("E40").Value = '  +9.98%  '
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("D42").Value = '0.784'
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("E43").Value = '  +2.02%  '
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("E45").Value = '  +4.68%  '
$ws.Range("D46").Value = '2.506.28'
$ws.Range("E46").Value = '  +2.25%  '
$ws.Range("D47").Value = '22.82'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '6.64'
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = '0.0263'
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '2.09'
$ws.Range("E51").Value = '  -1.13%  '
